$d = $word.ActiveDocument

# 1. Summary paragraph - simple text replacement
$d.Content.Find.Execute(
    "Highly motivated and results-oriented individual seeking a challenging role in the development of innovative mobile applications. Proven ability to collaborate effectively, solve complex problems, and contribute to a dynamic team environment. Eager to leverage skills in software development, project management, and communication to contribute to impactful projects.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "A motivated student with foundational knowledge in Kotlin, seeking an app development role. Leveraging academic learning and eagerness to grow into a professional position, contributing to innovative mobile application development.",
    2) | Out-Null

# 2. Experience paragraph - contains a w:br, rewrite whole paragraph text
$expPara = $d.Paragraphs(5)
$expRange = $expPara.Range
$expRange.End = $expRange.End - 1
$expRange.Text = "Developed and implemented an AI resume enhancer at Blue Silicon Infotech, yielding a 20% increase in resume completion rates. Optimized resume templates for enhanced readability and clarity, achieving a 15% improvement in resume accuracy through quantifiable results. Demonstrated expertise in AI-driven process optimization and template design, driving efficiency and effectiveness in resume development."

# 3. Education paragraph - simple text replacement
$d.Content.Find.Execute(
    "Bachelor of Engineering from AVIT. Graduated: 2026-05. GPA: 7.1.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Bachelor of Engineering, AVIT, May 2026, GPA 7.1.",
    2) | Out-Null

# 4. Skills paragraph - simple text replacement (whole huge block -> "C#")
$skillsPara = $d.Paragraphs(9)
$skillsRange = $skillsPara.Range
$skillsRange.End = $skillsRange.End - 1
$skillsRange.Text = "C#"

# 5. Projects paragraph - contains many w:br, rewrite whole paragraph text
$projPara = $d.Paragraphs(11)
$projRange = $projPara.Range
$projRange.End = $projRange.End - 1
$projRange.Text = "Developed an Enhanced QR Scanner and Generator project, significantly improving efficiency and accuracy. Implemented a novel algorithm and real-time data integration using QR scanner and generator, Prediction Pro, Simple Purchase Order Manager, and PDF Maker, resulting in 20% reduced processing time, improved accuracy, and enhanced real-time data integration, ultimately driving increased sales and lower operational costs."

Write-Output "Done"
